$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.193.92"
$ws.Range("E2").Value = "  -2.87%  "

# Row 3
$ws.Range("D3").Value = "3.295.18"
$ws.Range("E3").Value = "  -3.69%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "556.54"
$ws.Range("E5").Value = "  -3.87%  "

# Row 6
$ws.Range("D6").Value = "140.82"
$ws.Range("E6").Value = "  -8.58%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "3.295.39"
$ws.Range("E8").Value = "  -3.71%  "

# Row 9
$ws.Range("D9").Value = "0.467"
$ws.Range("E9").Value = "  -3.72%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.90"
$ws.Range("E10").Value = "  -1.56%  "

# Row 12
$ws.Range("E12").Value = "  -2.60%  "

# Row 13
$ws.Range("D13").Value = "3.862.43"
$ws.Range("E13").Value = "  -3.70%  "

# Row 14
$ws.Range("E14").Value = "  -0.14%  "

# Row 15
$ws.Range("D15").Value = "26.63"
$ws.Range("E15").Value = "  -5.99%  "

# Row 16
$ws.Range("D16").Value = "3.298.69"
$ws.Range("E16").Value = "  -3.12%  "

# Row 17
$ws.Range("E17").Value = "  -5.23%  "

# Row 18
$ws.Range("D18").Value = "60.214.36"
$ws.Range("E18").Value = "  -2.92%  "

# Row 19
$ws.Range("D19").Value = "6.04"
$ws.Range("E19").Value = "  -8.09%  "

# Row 20
$ws.Range("D20").Value = "13.67"
$ws.Range("E20").Value = "  -5.42%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.50"
$ws.Range("E21").Value = "  -5.04%  "

# Row 22
$ws.Range("D22").Value = "372.98"
$ws.Range("E22").Value = "  -2.54%  "

# Row 23
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "72.69"
$ws.Range("E23").Value = "  -4.28%  "

# Row 24
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.06%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.530"
$ws.Range("E25").Value = "  -6.99%  "

# Row 26
$ws.Range("D26").Value = "3.429.99"
$ws.Range("E26").Value = "  -3.69%  "

# Row 27
$ws.Range("E27").Value = "  -9.70%  "

# Row 28
$ws.Range("D28").Value = "0.174"
$ws.Range("E28").Value = "  -2.78%  "

# Row 29
$ws.Range("E29").Value = "  +0.00%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("E30").Value = "  -8.22%  "

# Row 31
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("D32").Value = "2.01"
$ws.Range("E32").Value = "  -5.22%  "

# Row 33
$ws.Range("D33").Value = "7.42"
$ws.Range("E33").Value = "  -5.84%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.50"
$ws.Range("E34").Value = "  -3.42%  "

# Row 35
$ws.Range("D35").Value = "1.22"
$ws.Range("E35").Value = "  -7.74%  "

# Row 36
$ws.Range("D36").Value = "165.87"
$ws.Range("E36").Value = "  -1.36%  "

# Row 37
$ws.Range("D37").Value = "5.02"
$ws.Range("E37").Value = "  -9.89%  "

# Row 38
$ws.Range("D38").Value = "1.52"
$ws.Range("E38").Value = "  -5.04%  "

# Row 39
$ws.Range("D39").Value = "6.61"
$ws.Range("E39").Value = "  -5.03%  "

# Row 40
$ws.Range("D40").Value = "3.329.45"
$ws.Range("E40").Value = "  -3.73%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0720"
$ws.Range("E41").Value = "  -8.49%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "25.40"
$ws.Range("E42").Value = "  -17.94%  "

# Row 43
$ws.Range("D43").Value = "41.61"
$ws.Range("E43").Value = "  -2.56%  "

# Row 44
$ws.Range("D44").Value = "0.747"
$ws.Range("E44").Value = "  -4.24%  "

# Row 45
$ws.Range("E45").Value = "  -4.50%  "

# Row 46
$ws.Range("D46").Value = "4.07"
$ws.Range("E46").Value = "  -8.02%  "

# Row 47
$ws.Range("E47").Value = "  -6.96%  "

# Row 48
$ws.Range("E48").Value = "  -0.03%  "

# Row 49
$ws.Range("D49").Value = "2.316.38"
$ws.Range("E49").Value = "  -9.38%  "

# Row 50
$ws.Range("D50").Value = "21.44"
$ws.Range("E50").Value = "  -7.10%  "

# Row 51
$ws.Range("D51").Value = "6.32"
$ws.Range("E51").Value = "  -7.22%  "
